$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new "line" entries (line7, line8) were added to the dataset. They are
# inserted right after the last existing "line" row (row 7) conceptually,
# which - because every row below shifts down by two - lands the two new
# rows at the bottom of the sheet (rows 16 & 17) while every row from 8
# downward receives the data that now belongs to it after the shift.

# Make room for the two new rows at the end of the table and copy the
# formatting (border/alignment) of the last data row onto them.
$ws.Range("A16:A17").EntireRow.Insert()
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)

# name | from_bus | to_bus | in_service
$data = @(
  @(8,  6,  "line7", 14, 11, $true),
  @(9,  7,  "line8", 16, 9,  $true),
  @(10, 8,  "extr1", 5,  12, $false),
  @(11, 9,  "extr2", 5,  9,  $false),
  @(12, 10, "extr3", 10, 11, $false),
  @(13, 11, "extr4", 7,  8,  $false),
  @(14, 12, "extr5", 9,  11, $false),
  @(15, 13, "extr6", 7,  11, $false),
  @(16, 14, "extr7", 5,  7,  $false),
  @(17, 15, "extr8", 8,  5,  $true)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

Write-Host "done"
